$d = $word.ActiveDocument

$pairs = @(
    @{old = "66×29=1914"; new = "68×43=2924"},
    @{old = "58×96=5568"; new = "76×51=3876"},
    @{old = "46×75=3450"; new = "57×40=2280"},
    @{old = "43×84=3612"; new = "86×15=1290"},
    @{old = "17×37=629";  new = "22×42=924"},
    @{old = "48×53=2544"; new = "38×33=1254"},
    @{old = "65×74=4810"; new = "77×14=1078"},
    @{old = "27×19=513";  new = "30×11=330"},
    @{old = "43×38=1634"; new = "55×93=5115"},
    @{old = "36×42=1512"; new = "74×19=1406"},
    @{old = "35×18=630";  new = "25×36=900"},
    @{old = "86×31=2666"; new = "35×70=2450"},
    @{old = "54×63=3402"; new = "23×45=1035"},
    @{old = "12×93=1116"; new = "33×73=2409"},
    @{old = "66×50=3300"; new = "61×19=1159"},
    @{old = "63×67=4221"; new = "53×75=3975"},
    @{old = "62×66=4092"; new = "26×12=312"},
    @{old = "88×41=3608"; new = "82×75=6150"},
    @{old = "83×49=4067"; new = "98×38=3724"},
    @{old = "13×34=442";  new = "32×82=2624"},
    @{old = "60×63=3780"; new = "88×45=3960"},
    @{old = "17×77=1309"; new = "93×18=1674"},
    @{old = "13×48=624";  new = "74×63=4662"},
    @{old = "24×48=1152"; new = "89×85=7565"},
    @{old = "17×51=867";  new = "51×54=2754"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}
